# Generate Report for Handback
# Refresh the two source UUIDs / handoff+handback xlf hashes / timestamps
# that the handback-status report tracks, across all three sheets
# (Overview, zh-cn, de-de), keeping each hyperlink's underlying target
# URL intact but updating the displayed text to match the new values.

$wb = $excel.ActiveWorkbook

# ---- new data ----------------------------------------------------------
$newUuid1 = "490b1726-4faf-4619-b87d-251fc04d19b3"
$newUuid2 = "ffffd56bfdba-86b2-4f81-a24c-06c5fba419e5"
$newHash  = "5c58728acdedc51e4ff23cc4322e5fae940ad35a"

$zhTime1Off  = "2016-03-20 10:49:56"
$zhTime1Back = "2016-03-20 10:50:15"
$deTime1Off  = "2016-03-20 10:49:59"
$deTime1Back = "2016-03-20 10:50:20"

$newUuid1Md = "$newUuid1.md"
$newUuid2Md = "$newUuid2.md"
$newZhXlf   = "$newUuid1.$newHash.zh-cn.xlf"
$newDeXlf   = "$newUuid1.$newHash.de-de.xlf"

# ---- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newUuid1Md
$wsOverview.Range("A3").Value = $newUuid2Md

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", $newUuid1Md) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", $newUuid2Md) | Out-Null

# ---- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newUuid1Md
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $zhTime1Off
$wsZh.Range("F2").Value = $newUuid1Md
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $zhTime1Back

$wsZh.Range("A3").Value = $newUuid2Md
$wsZh.Range("D3").Value = $newZhXlf
$wsZh.Range("E3").Value = $zhTime1Off
$wsZh.Range("F3").Value = $newUuid2Md
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $zhTime1Back

$wsZh.Range("A2").Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", $newUuid1Md) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2260a60e2799454237861fb46b4fd2470a45ff2a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.zh-cn.xlf", "", "", $newZhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a13d7494d381be89aa0111cba84195137e88d49e/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", $newUuid1Md) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0af3984e50313562de923dc6896c674ebed9dc4f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.zh-cn.xlf", "", "", $newZhXlf) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", $newUuid2Md) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2260a60e2799454237861fb46b4fd2470a45ff2a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.zh-cn.xlf", "", "", $newZhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a13d7494d381be89aa0111cba84195137e88d49e/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", $newUuid2Md) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0af3984e50313562de923dc6896c674ebed9dc4f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.zh-cn.xlf", "", "", $newZhXlf) | Out-Null

# ---- de-de sheet -----------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newUuid1Md
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $deTime1Off
$wsDe.Range("F2").Value = $newUuid1Md
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $deTime1Back

$wsDe.Range("A3").Value = $newUuid2Md
$wsDe.Range("D3").Value = $newDeXlf
$wsDe.Range("E3").Value = $deTime1Off
$wsDe.Range("F3").Value = $newUuid2Md
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $deTime1Back

$wsDe.Range("A2").Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", $newUuid1Md) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5e4094ed9d048769903debf9fc1f9097c5a43b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.de-de.xlf", "", "", $newDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6d7c9a2c2ec53acd4ea73f28d3bc6126d50db0b2/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", $newUuid1Md) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/10fed944dd353e10ceb53f5bb66f3261aa1f0559/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.de-de.xlf", "", "", $newDeXlf) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", $newUuid2Md) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5e4094ed9d048769903debf9fc1f9097c5a43b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.de-de.xlf", "", "", $newDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6d7c9a2c2ec53acd4ea73f28d3bc6126d50db0b2/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", $newUuid2Md) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/10fed944dd353e10ceb53f5bb66f3261aa1f0559/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.de-de.xlf", "", "", $newDeXlf) | Out-Null

Write-Output "handback status refreshed"
